# Update weights and eggs
# Adds three new daily-log rows (35-37) to Sheet1, continuing the
# existing date sequence, and leaves the selection on the cell the
# author ended up on after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-number-format from the last existing row (A34) down
# into the three new rows so the new dates render the same way
# (reuses the workbook's existing style instead of creating a new one).
$ws.Range("A34").Copy()
$ws.Range("A35:A37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 35: 2021-02-22
$ws.Range("A35").Value = 44249
$ws.Range("B35").Value = 48
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0

# Row 36: 2021-02-23
$ws.Range("A36").Value = 44250
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0

# Row 37: 2021-02-24
$ws.Range("A37").Value = 44251
$ws.Range("B37").Value = 48
$ws.Range("C37").Value = 15
$ws.Range("D37").Value = 12
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0

# Leave the view scrolled back to the top of the sheet with the
# selection where the author left it.
$ws.Range("J32").Select()
